$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update input values (formulas downstream recalc automatically)
$ws.Range("B17").Value = 4
$ws.Range("C19").Value = 2

# Column A width change (stored width 45 == ColumnWidth ~44.14 after Excel's pixel rounding)
$ws.Columns.Item(1).ColumnWidth = 44.14

# Update the active selection to C26
$ws.Range("C26").Select()
